# Applies the crypto price/volume update for Mon Aug 19 15:39:35 UTC 2024 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for each changed cell (column D/E hold numeric-looking text that must
# stay text, so their number format is forced to "@" before assignment).
$updates = [ordered]@{
    'D2' = '58.516.82'
    'E2' = '  -2.16%  '
    'D3' = '2.580.10'
    'E3' = '  -3.18%  '
    'E4' = '  +0.04%  '
    'D5' = '542.01'
    'E5' = '  +0.70%  '
    'D6' = '143.89'
    'E6' = '  -1.22%  '
    'E7' = '  -0.04%  '
    'E8' = '  +1.63%  '
    'D9' = '6.75'
    'E9' = '  +1.09%  '
    'E10' = '  -2.95%  '
    'E11' = '  +3.50%  '
    'E12' = '  -1.91%  '
    'D13' = '3.033.29'
    'E13' = '  -3.15%  '
    'D14' = '58.444.07'
    'E14' = '  -2.14%  '
    'D15' = '20.58'
    'E15' = '  -2.95%  '
    'B16' = 'WrappedEther'
    'C16' = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
    'D16' = '2.595.09'
    'E16' = '  -3.11%  '
    'B17' = 'ShibaInu'
    'C17' = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
    'D17' = '0.0000131'
    'E17' = '  -2.80%  '
    'D18' = '4.45'
    'E18' = '  +0.67%  '
    'D19' = '334.26'
    'E19' = '  -3.26%  '
    'D20' = '10.05'
    'E20' = '  -3.49%  '
    'D21' = '6.08'
    'E21' = '  -4.26%  '
    'E22' = '  -0.11%  '
    'D23' = '66.32'
    'E23' = '  -0.43%  '
    'D24' = '0.423'
    'E24' = '  +1.58%  '
    'D25' = '0.999'
    'E25' = '  -0.18%  '
    'E26' = '  -4.99%  '
    'D27' = '7.07'
    'E27' = '  -3.32%  '
    'D28' = '0.0₃0741'
    'E28' = '  -1.69%  '
    'D29' = '0.998'
    'E29' = '  -0.03%  '
    'E30' = '  -1.33%  '
    'D31' = '5.97'
    'E31' = '  +2.10%  '
    'E32' = '  +1.88%  '
    'D33' = '18.91'
    'E33' = '  -0.69%  '
    'E34' = '  -2.99%  '
    'D35' = '0.846'
    'E35' = '  +2.33%  '
    'E36' = '  -4.73%  '
    'D37' = '0.820'
    'E37' = '  -2.98%  '
    'E38' = '  -3.02%  '
    'E39' = '  -0.77%  '
    'D40' = '278.50'
    'E40' = '  -4.95%  '
    'E41' = '  -0.19%  '
    'E42' = '  -2.41%  '
    'D43' = '10.62'
    'E43' = '  -0.94%  '
    'D44' = '0.0942'
    'E44' = '  -0.85%  '
    'E45' = '  -2.47%  '
    'D46' = '18.51'
    'E46' = '  -5.09%  '
    'D47' = '0.0227'
    'E47' = '  +0.20%  '
    'D48' = '1.902.36'
    'E48' = '  -3.84%  '
    'D49' = '17.86'
    'E49' = '  -2.92%  '
    'E50' = '  -2.99%  '
    'D51' = '109.13'
    'E51' = '  -1.12%  '
}

foreach ($cellRef in $updates.Keys) {
    $col = $cellRef -replace '[0-9]+$', ''
    $range = $ws.Range($cellRef)
    if ($col -eq 'D' -or $col -eq 'E') {
        $range.NumberFormat = "@"
    }
    $range.Value = $updates[$cellRef]
}
